# Add two new columns (I: "I0", J: "IF") to the sheet, matching the
# existing header style used by column H, and fill in the corresponding
# data values for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, centered, bordered) used by the
# other header cells (e.g. H1) by copying its format onto the new cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 6
